$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# xlPasteFormats constant
$xlPasteFormats = -4122

# --- Row 66: Modifica Poi / 3 / completato (E66 new), matches style of row 65 (C,D,E) ---
$ws.Range("C65").Copy() | Out-Null
$ws.Range("C66").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C66").Value = "Modifica Poi"

$ws.Range("D65").Copy() | Out-Null
$ws.Range("D66").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D66").Value = 3

$ws.Range("E65").Copy() | Out-Null
$ws.Range("E66").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E66").Value = "completato"

# --- Row 68: Segnala contenuto / 5 / (F68 new) ---
$ws.Range("C65").Copy() | Out-Null
$ws.Range("C68").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C68").Value = "Segnala contenuto"
$ws.Range("D68").Value = 5

# --- Row 69: Modifica Comune / 6 / (F69 new, reuses "dettaglio con diagramma ") ---
$ws.Range("C65").Copy() | Out-Null
$ws.Range("C69").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C69").Value = "Modifica Comune"
$ws.Range("D69").Value = 6

# String "dettaglio con diagramma " must be created now (3rd new unique string), used on F68 first
$ws.Range("F68").Value = "dettaglio con diagramma "
$ws.Range("F69").Value = "dettaglio con diagramma "

# --- Row 72: new header row (date 15/02/2025, iteration 5) + Associa Poi (existing) + D72=1 ---
$ws.Range("A63").Copy() | Out-Null
$ws.Range("A72").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A72").Value = "2/15/2025"

$ws.Range("B63").Copy() | Out-Null
$ws.Range("B72").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B72").Value = 5

$ws.Range("C72").Value = "Associa Poi"
$ws.Range("D72").Value = 1

# --- Row 73 (new): Modifica Attività (reuses existing string) / 2 ---
$ws.Range("C73").Value = "Modifica Attività "
$ws.Range("D73").Value = 2

# --- Row 74 (new): Rimuovi Attività (reuses existing string) / 3 ---
$ws.Range("C74").Value = "Rimuovi Attività"
$ws.Range("D74").Value = 3

# --- Row 75 (new): Modifica dati Utente (new string #4), style matches row65 (C,D), F reuses "dettaglio con diagramma " ---
$ws.Range("C65").Copy() | Out-Null
$ws.Range("C75").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C75").Value = "Modifica dati Utente"

$ws.Range("D65").Copy() | Out-Null
$ws.Range("D75").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D75").Value = 4

$ws.Range("F75").Value = "dettaglio con diagramma "

# --- Row 76 (new): Salva informazioni per visite future (new string #5) / 5 ---
$ws.Range("C76").Value = "Salva informazioni per visite future"
$ws.Range("D76").Value = 5

# --- Column C width: stored OOXML width 32 (ColumnWidth API has +5/6 padding offset) ---
$ws.Columns.Item(3).ColumnWidth = 31.16666666666667

# --- View: selection on E74 ---
$ws.Activate()
$ws.Range("E74").Select()
